# Auto-generated edit script updating Leve profit calculation columns (H:N)
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
# Row 6 (Item ID G6=4564)
$ws.Cells.Item(6, 8).Value = 1176.2858  # H6
$ws.Cells.Item(6, 9).Value = 251.66667  # I6
$ws.Cells.Item(6, 10).Value = 1869.75  # J6
$ws.Cells.Item(6, 11).Value = 755.00001  # K6
$ws.Cells.Item(6, 12).Value = 5609.25  # L6
$ws.Cells.Item(6, 13).Value = -643.00001  # M6
$ws.Cells.Item(6, 14).Value = -5833.25  # N6

# Row 32 (Item ID G32=5484)
$ws.Cells.Item(32, 8).Value = 23086584  # H32
$ws.Cells.Item(32, 9).Value = 28579942  # I32
$ws.Cells.Item(32, 10).Value = 16677665  # J32
$ws.Cells.Item(32, 11).Value = 28579942  # K32
$ws.Cells.Item(32, 12).Value = 16677665  # L32
$ws.Cells.Item(32, 13).Value = -28579616  # M32
$ws.Cells.Item(32, 14).Value = -16678317  # N32

# Row 33 (Item ID G33=5512)
$ws.Cells.Item(33, 8).Value = 1021.13336  # H33
$ws.Cells.Item(33, 9).Value = 1079.7858  # I33
$ws.Cells.Item(33, 10).Value = 200  # J33
$ws.Cells.Item(33, 11).Value = 1079.7858  # K33
$ws.Cells.Item(33, 12).Value = 200  # L33
$ws.Cells.Item(33, 13).Value = -850.7858000000001  # M33
$ws.Cells.Item(33, 14).Value = -658  # N33

# Row 38 (Item ID G38=4599)
$ws.Cells.Item(38, 8).Value = 195.6  # H38
$ws.Cells.Item(38, 9).Value = 195.6  # I38
$ws.Cells.Item(38, 11).Value = 586.8  # K38
$ws.Cells.Item(38, 13).Value = -214.8  # M38

# Row 41 (Item ID G41=5478)
$ws.Cells.Item(41, 8).Value = 312.5  # H41
$ws.Cells.Item(41, 9).Value = 264.3  # I41
$ws.Cells.Item(41, 11).Value = 264.3  # K41
$ws.Cells.Item(41, 13).Value = 175.7  # M41

# Row 42 (Item ID G42=4600)
$ws.Cells.Item(42, 8).Value = 957.625  # H42
$ws.Cells.Item(42, 9).Value = 1321.4  # I42
$ws.Cells.Item(42, 10).Value = 351.33334  # J42
$ws.Cells.Item(42, 11).Value = 3964.2  # K42
$ws.Cells.Item(42, 12).Value = 1054.00002  # L42
$ws.Cells.Item(42, 13).Value = -3734.2  # M42
$ws.Cells.Item(42, 14).Value = -1514.00002  # N42

# Row 53 (Item ID G53=5479)
$ws.Cells.Item(53, 8).Value = 735.625  # H53
$ws.Cells.Item(53, 10).Value = 197.85715  # J53
$ws.Cells.Item(53, 12).Value = 197.85715  # L53
$ws.Cells.Item(53, 14).Value = -1471.85715  # N53

# Row 62 (Item ID G62=27781)
$ws.Cells.Item(62, 8).Value = 4644.091  # H62
$ws.Cells.Item(62, 9).Value = 4122.3335  # I62
$ws.Cells.Item(62, 10).Value = 5270.2  # J62
$ws.Cells.Item(62, 11).Value = 4122.3335  # K62
$ws.Cells.Item(62, 12).Value = 5270.2  # L62
$ws.Cells.Item(62, 13).Value = -3498.3335  # M62
$ws.Cells.Item(62, 14).Value = -6518.2  # N62

# Row 65 (Item ID G65=27781)
$ws.Cells.Item(65, 8).Value = 4644.091  # H65
$ws.Cells.Item(65, 9).Value = 4122.3335  # I65
$ws.Cells.Item(65, 10).Value = 5270.2  # J65
$ws.Cells.Item(65, 11).Value = 20611.6675  # K65
$ws.Cells.Item(65, 12).Value = 26351  # L65
$ws.Cells.Item(65, 13).Value = -17491.6675  # M65
$ws.Cells.Item(65, 14).Value = -32591  # N65

# Row 80 (Item ID G80=12605)
$ws.Cells.Item(80, 8).Value = 705.931  # H80
$ws.Cells.Item(80, 9).Value = 429.72  # I80
$ws.Cells.Item(80, 10).Value = 2432.25  # J80
$ws.Cells.Item(80, 11).Value = 1289.16  # K80
$ws.Cells.Item(80, 12).Value = 7296.75  # L80
$ws.Cells.Item(80, 13).Value = -291.1600000000001  # M80
$ws.Cells.Item(80, 14).Value = -9292.75  # N80

# Row 82 (Item ID G82=12623)
$ws.Cells.Item(82, 8).Value = 1770.1428  # H82
$ws.Cells.Item(82, 9).Value = 1770.1428  # I82
$ws.Cells.Item(82, 11).Value = 5310.428400000001  # K82
$ws.Cells.Item(82, 13).Value = -4904.428400000001  # M82

# Row 83 (Item ID G83=12605)
$ws.Cells.Item(83, 8).Value = 705.931  # H83
$ws.Cells.Item(83, 9).Value = 429.72  # I83
$ws.Cells.Item(83, 10).Value = 2432.25  # J83
$ws.Cells.Item(83, 11).Value = 3867.48  # K83
$ws.Cells.Item(83, 12).Value = 21890.25  # L83
$ws.Cells.Item(83, 13).Value = 1124.52  # M83
$ws.Cells.Item(83, 14).Value = -31874.25  # N83

# Row 85 (Item ID G85=12623)
$ws.Cells.Item(85, 8).Value = 1770.1428  # H85
$ws.Cells.Item(85, 9).Value = 1770.1428  # I85
$ws.Cells.Item(85, 11).Value = 5310.428400000001  # K85
$ws.Cells.Item(85, 13).Value = -3906.428400000001  # M85

# Row 93 (Item ID G93=18043)
$ws.Cells.Item(93, 8).Value = 0  # H93
$ws.Cells.Item(93, 10).Value = 0  # J93
$ws.Cells.Item(93, 12).Value = 0  # L93
$ws.Cells.Item(93, 14).ClearContents()  # N93

# Row 94 (Item ID G94=19905)
$ws.Cells.Item(94, 8).Value = 11907144  # H94
$ws.Cells.Item(94, 9).Value = 11907144  # I94
$ws.Cells.Item(94, 10).Value = 0  # J94
$ws.Cells.Item(94, 11).Value = 11907144  # K94
$ws.Cells.Item(94, 12).Value = 0  # L94
$ws.Cells.Item(94, 13).Value = -11906693  # M94
$ws.Cells.Item(94, 14).ClearContents()  # N94

# Row 98 (Item ID G98=36237)
$ws.Cells.Item(98, 8).Value = 1064.0588  # H98
$ws.Cells.Item(98, 9).Value = 1093.0625  # I98
$ws.Cells.Item(98, 11).Value = 1093.0625  # K98
$ws.Cells.Item(98, 13).Value = 404.9375  # M98

# Row 104 (Item ID G104=24263)
$ws.Cells.Item(104, 8).Value = 856.8  # H104
$ws.Cells.Item(104, 9).Value = 913.55554  # I104
$ws.Cells.Item(104, 11).Value = 2740.66662  # K104
$ws.Cells.Item(104, 13).Value = -993.66662  # M104

# Row 111 (Item ID G111=27768)
$ws.Cells.Item(111, 8).Value = 1421.4615  # H111
$ws.Cells.Item(111, 9).Value = 1515.2727  # I111
$ws.Cells.Item(111, 10).Value = 905.5  # J111
$ws.Cells.Item(111, 11).Value = 4545.8181  # K111
$ws.Cells.Item(111, 12).Value = 2716.5  # L111
$ws.Cells.Item(111, 13).Value = -1478.8181  # M111
$ws.Cells.Item(111, 14).Value = -8850.5  # N111

# Row 113 (Item ID G113=27775)
$ws.Cells.Item(113, 8).Value = 15411.667  # H113
$ws.Cells.Item(113, 9).Value = 6241.2  # I113
$ws.Cells.Item(113, 11).Value = 6241.2  # K113
$ws.Cells.Item(113, 13).Value = -2987.2  # M113

# Row 122 (Item ID G122=36237)
$ws.Cells.Item(122, 8).Value = 1064.0588  # H122
$ws.Cells.Item(122, 9).Value = 1093.0625  # I122
$ws.Cells.Item(122, 11).Value = 3279.1875  # K122
$ws.Cells.Item(122, 13).Value = -829.1875  # M122

# Row 132 (Item ID G132=44049)
$ws.Cells.Item(132, 8).Value = 11926.19  # H132
$ws.Cells.Item(132, 9).Value = 4792.467  # I132
$ws.Cells.Item(132, 10).Value = 29760.5  # J132
$ws.Cells.Item(132, 11).Value = 14377.401  # K132
$ws.Cells.Item(132, 12).Value = 89281.5  # L132
$ws.Cells.Item(132, 13).Value = -11847.401  # M132
$ws.Cells.Item(132, 14).Value = -94341.5  # N132

# Row 137 (Item ID G137=44013)
$ws.Cells.Item(137, 8).Value = 1908.9524  # H137
$ws.Cells.Item(137, 9).Value = 1787.5555  # I137
$ws.Cells.Item(137, 10).Value = 2000  # J137
$ws.Cells.Item(137, 11).Value = 5362.666499999999  # K137
$ws.Cells.Item(137, 12).Value = 6000  # L137
$ws.Cells.Item(137, 13).Value = -2812.666499999999  # M137
$ws.Cells.Item(137, 14).Value = -11100  # N137

# Row 138 (Item ID G138=44169)
$ws.Cells.Item(138, 8).Value = 1955.22  # H138
$ws.Cells.Item(138, 9).Value = 1600.591  # I138
$ws.Cells.Item(138, 11).Value = 4801.772999999999  # K138
$ws.Cells.Item(138, 13).Value = 338.2270000000008  # M138

# Row 141 (Item ID G141=44161)
$ws.Cells.Item(141, 8).Value = 6367.564  # H141
$ws.Cells.Item(141, 9).Value = 7958.8237  # I141
$ws.Cells.Item(141, 11).Value = 23876.4711  # K141
$ws.Cells.Item(141, 13).Value = -18696.4711  # M141

$ws = $wb.Worksheets.Item(2)  # ARM
# Row 32 (Item ID G32=44147)
$ws.Cells.Item(32, 8).Value = 8211.4  # H32
$ws.Cells.Item(32, 9).Value = 5993.8887  # I32
$ws.Cells.Item(32, 10).Value = 28169  # J32
$ws.Cells.Item(32, 11).Value = 5993.8887  # K32
$ws.Cells.Item(32, 12).Value = 28169  # L32
$ws.Cells.Item(32, 13).Value = -5706.8887  # M32
$ws.Cells.Item(32, 14).Value = -28743  # N32

# Row 45 (Item ID G45=27714)
$ws.Cells.Item(45, 8).Value = 12288.692  # H45
$ws.Cells.Item(45, 9).Value = 16455.334  # I45
$ws.Cells.Item(45, 11).Value = 16455.334  # K45
$ws.Cells.Item(45, 13).Value = -16078.334  # M45

# Row 61 (Item ID G61=43999)
$ws.Cells.Item(61, 8).Value = 19295.334  # H61
$ws.Cells.Item(61, 9).Value = 18591.666  # I61
$ws.Cells.Item(61, 10).Value = 19999  # J61
$ws.Cells.Item(61, 11).Value = 18591.666  # K61
$ws.Cells.Item(61, 12).Value = 19999  # L61
$ws.Cells.Item(61, 13).Value = -18379.666  # M61
$ws.Cells.Item(61, 14).Value = -20423  # N61

# Row 97 (Item ID G97=19941)
$ws.Cells.Item(97, 8).Value = 34518824  # H97
$ws.Cells.Item(97, 9).Value = 41668516  # I97
$ws.Cells.Item(97, 10).Value = 200290  # J97
$ws.Cells.Item(97, 11).Value = 41668516  # K97
$ws.Cells.Item(97, 12).Value = 200290  # L97
$ws.Cells.Item(97, 13).Value = -41668020  # M97
$ws.Cells.Item(97, 14).Value = -201282  # N97

# Row 102 (Item ID G102=19945)
$ws.Cells.Item(102, 8).Value = 3388.2666  # H102
$ws.Cells.Item(102, 9).Value = 3416.0715  # I102
$ws.Cells.Item(102, 11).Value = 3416.0715  # K102
$ws.Cells.Item(102, 13).Value = -1794.0715  # M102

# Row 110 (Item ID G110=27708)
$ws.Cells.Item(110, 8).Value = 4814.8  # H110
$ws.Cells.Item(110, 9).Value = 4366  # I110
$ws.Cells.Item(110, 10).Value = 5488  # J110
$ws.Cells.Item(110, 11).Value = 4366  # K110
$ws.Cells.Item(110, 12).Value = 5488  # L110
$ws.Cells.Item(110, 13).Value = -2321  # M110
$ws.Cells.Item(110, 14).Value = -9578  # N110

# Row 122 (Item ID G122=36168)
$ws.Cells.Item(122, 8).Value = 4799.3  # H122
$ws.Cells.Item(122, 9).Value = 3499.25  # I122
$ws.Cells.Item(122, 11).Value = 10497.75  # K122
$ws.Cells.Item(122, 13).Value = -8047.75  # M122

# Row 132 (Item ID G132=43997)
$ws.Cells.Item(132, 8).Value = 3360.9062  # H132
$ws.Cells.Item(132, 9).Value = 3385.484  # I132
$ws.Cells.Item(132, 11).Value = 10156.452  # K132
$ws.Cells.Item(132, 13).Value = -7626.451999999999  # M132

# Row 136 (Item ID G136=43999)
$ws.Cells.Item(136, 8).Value = 19295.334  # H136
$ws.Cells.Item(136, 9).Value = 18591.666  # I136
$ws.Cells.Item(136, 10).Value = 19999  # J136
$ws.Cells.Item(136, 11).Value = 55774.99800000001  # K136
$ws.Cells.Item(136, 12).Value = 59997  # L136
$ws.Cells.Item(136, 13).Value = -53224.99800000001  # M136
$ws.Cells.Item(136, 14).Value = -65097  # N136

$ws = $wb.Worksheets.Item(3)  # BSM
# Row 20 (Item ID G20=14149)
$ws.Cells.Item(20, 8).Value = 3921.4814  # H20
$ws.Cells.Item(20, 9).Value = 4154.222  # I20
$ws.Cells.Item(20, 10).Value = 3456  # J20
$ws.Cells.Item(20, 11).Value = 4154.222  # K20
$ws.Cells.Item(20, 12).Value = 3456  # L20
$ws.Cells.Item(20, 13).Value = -3907.222  # M20
$ws.Cells.Item(20, 14).Value = -3950  # N20

# Row 22 (Item ID G22=5092)
$ws.Cells.Item(22, 8).Value = 496.41666  # H22
$ws.Cells.Item(22, 9).Value = 496.41666  # I22
$ws.Cells.Item(22, 11).Value = 496.41666  # K22
$ws.Cells.Item(22, 13).Value = -323.41666  # M22

# Row 86 (Item ID G86=12526)
$ws.Cells.Item(86, 8).Value = 2473.75  # H86
$ws.Cells.Item(86, 9).Value = 2498.1  # I86
$ws.Cells.Item(86, 10).Value = 2352  # J86
$ws.Cells.Item(86, 11).Value = 2498.1  # K86
$ws.Cells.Item(86, 12).Value = 2352  # L86
$ws.Cells.Item(86, 13).Value = -1375.1  # M86
$ws.Cells.Item(86, 14).Value = -4598  # N86

# Row 89 (Item ID G89=12526)
$ws.Cells.Item(89, 8).Value = 2473.75  # H89
$ws.Cells.Item(89, 9).Value = 2498.1  # I89
$ws.Cells.Item(89, 10).Value = 2352  # J89
$ws.Cells.Item(89, 11).Value = 12490.5  # K89
$ws.Cells.Item(89, 12).Value = 11760  # L89
$ws.Cells.Item(89, 13).Value = -6874.5  # M89
$ws.Cells.Item(89, 14).Value = -22992  # N89

# Row 94 (Item ID G94=19939)
$ws.Cells.Item(94, 8).Value = 1343.3529  # H94
$ws.Cells.Item(94, 9).Value = 1103.909  # I94
$ws.Cells.Item(94, 10).Value = 1782.3334  # J94
$ws.Cells.Item(94, 11).Value = 1103.909  # K94
$ws.Cells.Item(94, 12).Value = 1782.3334  # L94
$ws.Cells.Item(94, 13).Value = -652.9090000000001  # M94
$ws.Cells.Item(94, 14).Value = -2684.3334  # N94

# Row 105 (Item ID G105=19947)
$ws.Cells.Item(105, 8).Value = 4533.9355  # H105
$ws.Cells.Item(105, 9).Value = 3347.8928  # I105
$ws.Cells.Item(105, 11).Value = 3347.8928  # K105
$ws.Cells.Item(105, 13).Value = -1600.8928  # M105

# Row 107 (Item ID G107=27706)
$ws.Cells.Item(107, 8).Value = 7154.478  # H107
$ws.Cells.Item(107, 9).Value = 7174.2  # I107
$ws.Cells.Item(107, 10).Value = 7117.5  # J107
$ws.Cells.Item(107, 11).Value = 7174.2  # K107
$ws.Cells.Item(107, 12).Value = 7117.5  # L107
$ws.Cells.Item(107, 13).Value = -5254.2  # M107
$ws.Cells.Item(107, 14).Value = -10957.5  # N107

# Row 129 (Item ID G129=35382)
$ws.Cells.Item(129, 8).Value = 0  # H129
$ws.Cells.Item(129, 10).Value = 0  # J129
$ws.Cells.Item(129, 12).Value = 0  # L129
$ws.Cells.Item(129, 14).ClearContents()  # N129

# Row 134 (Item ID G134=43998)
$ws.Cells.Item(134, 8).Value = 6117  # H134
$ws.Cells.Item(134, 9).Value = 6402.56  # I134
$ws.Cells.Item(134, 10).Value = 2547.5  # J134
$ws.Cells.Item(134, 11).Value = 19207.68  # K134
$ws.Cells.Item(134, 12).Value = 7642.5  # L134
$ws.Cells.Item(134, 13).Value = -16672.68  # M134
$ws.Cells.Item(134, 14).Value = -12712.5  # N134

$ws = $wb.Worksheets.Item(4)  # CRP
# Row 16 (Item ID G16=27691)
$ws.Cells.Item(16, 8).Value = 5897.25  # H16
$ws.Cells.Item(16, 9).Value = 3420.75  # I16
$ws.Cells.Item(16, 10).Value = 8373.75  # J16
$ws.Cells.Item(16, 11).Value = 3420.75  # K16
$ws.Cells.Item(16, 12).Value = 8373.75  # L16
$ws.Cells.Item(16, 13).Value = -3133.75  # M16
$ws.Cells.Item(16, 14).Value = -8947.75  # N16

# Row 22 (Item ID G22=5367)
$ws.Cells.Item(22, 8).Value = 1930.1904  # H22
$ws.Cells.Item(22, 9).Value = 1825.2727  # I22
$ws.Cells.Item(22, 11).Value = 1825.2727  # K22
$ws.Cells.Item(22, 13).Value = -1475.2727  # M22

# Row 31 (Item ID G31=44023)
$ws.Cells.Item(31, 8).Value = 4534.2  # H31
$ws.Cells.Item(31, 9).Value = 4411.1665  # I31
$ws.Cells.Item(31, 10).Value = 4616.222  # J31
$ws.Cells.Item(31, 11).Value = 4411.1665  # K31
$ws.Cells.Item(31, 12).Value = 4616.222  # L31
$ws.Cells.Item(31, 13).Value = -4116.1665  # M31
$ws.Cells.Item(31, 14).Value = -5206.222  # N31

# Row 34 (Item ID G34=44023)
$ws.Cells.Item(34, 8).Value = 4534.2  # H34
$ws.Cells.Item(34, 9).Value = 4411.1665  # I34
$ws.Cells.Item(34, 10).Value = 4616.222  # J34
$ws.Cells.Item(34, 11).Value = 4411.1665  # K34
$ws.Cells.Item(34, 12).Value = 4616.222  # L34
$ws.Cells.Item(34, 13).Value = -4209.1665  # M34
$ws.Cells.Item(34, 14).Value = -5020.222  # N34

# Row 58 (Item ID G58=44021)
$ws.Cells.Item(58, 8).Value = 6416.231  # H58
$ws.Cells.Item(58, 9).Value = 5757.9414  # I58
$ws.Cells.Item(58, 10).Value = 7659.6665  # J58
$ws.Cells.Item(58, 11).Value = 5757.9414  # K58
$ws.Cells.Item(58, 12).Value = 7659.6665  # L58
$ws.Cells.Item(58, 13).Value = -5554.9414  # M58
$ws.Cells.Item(58, 14).Value = -8065.6665  # N58

# Row 105 (Item ID G105=19928)
$ws.Cells.Item(105, 8).Value = 1811.4138  # H105
$ws.Cells.Item(105, 9).Value = 1723.55  # I105
$ws.Cells.Item(105, 10).Value = 2006.6666  # J105
$ws.Cells.Item(105, 11).Value = 1723.55  # K105
$ws.Cells.Item(105, 12).Value = 2006.6666  # L105
$ws.Cells.Item(105, 13).Value = 23.45000000000005  # M105
$ws.Cells.Item(105, 14).Value = -5500.6666  # N105

# Row 107 (Item ID G107=27689)
$ws.Cells.Item(107, 8).Value = 1566.4286  # H107
$ws.Cells.Item(107, 9).Value = 1178.1111  # I107
$ws.Cells.Item(107, 10).Value = 2265.4  # J107
$ws.Cells.Item(107, 11).Value = 1178.1111  # K107
$ws.Cells.Item(107, 12).Value = 2265.4  # L107
$ws.Cells.Item(107, 13).Value = 741.8888999999999  # M107
$ws.Cells.Item(107, 14).Value = -6105.4  # N107

# Row 113 (Item ID G113=27691)
$ws.Cells.Item(113, 8).Value = 5897.25  # H113
$ws.Cells.Item(113, 9).Value = 3420.75  # I113
$ws.Cells.Item(113, 10).Value = 8373.75  # J113
$ws.Cells.Item(113, 11).Value = 3420.75  # K113
$ws.Cells.Item(113, 12).Value = 8373.75  # L113
$ws.Cells.Item(113, 13).Value = -1250.75  # M113
$ws.Cells.Item(113, 14).Value = -12713.75  # N113

# Row 122 (Item ID G122=36196)
$ws.Cells.Item(122, 8).Value = 3972.8  # H122
$ws.Cells.Item(122, 9).Value = 3661.6  # I122
$ws.Cells.Item(122, 10).Value = 4595.2  # J122
$ws.Cells.Item(122, 11).Value = 10984.8  # K122
$ws.Cells.Item(122, 12).Value = 13785.6  # L122
$ws.Cells.Item(122, 13).Value = -8534.799999999999  # M122
$ws.Cells.Item(122, 14).Value = -18685.6  # N122

# Row 132 (Item ID G132=44019)
$ws.Cells.Item(132, 8).Value = 4974.0884  # H132
$ws.Cells.Item(132, 9).Value = 4974.0884  # I132
$ws.Cells.Item(132, 11).Value = 14922.2652  # K132
$ws.Cells.Item(132, 13).Value = -12392.2652  # M132

# Row 134 (Item ID G134=44020)
$ws.Cells.Item(134, 8).Value = 100011  # H134
$ws.Cells.Item(134, 9).Value = 100011  # I134
$ws.Cells.Item(134, 11).Value = 300033  # K134
$ws.Cells.Item(134, 13).Value = -297498  # M134

# Row 136 (Item ID G136=44021)
$ws.Cells.Item(136, 8).Value = 6416.231  # H136
$ws.Cells.Item(136, 9).Value = 5757.9414  # I136
$ws.Cells.Item(136, 10).Value = 7659.6665  # J136
$ws.Cells.Item(136, 11).Value = 17273.8242  # K136
$ws.Cells.Item(136, 12).Value = 22978.9995  # L136
$ws.Cells.Item(136, 13).Value = -14723.8242  # M136
$ws.Cells.Item(136, 14).Value = -28078.9995  # N136

$ws = $wb.Worksheets.Item(5)  # CUL
# Row 5 (Item ID G5=43974)
$ws.Cells.Item(5, 8).Value = 645.8889  # H5
$ws.Cells.Item(5, 9).Value = 545.5714  # I5
$ws.Cells.Item(5, 11).Value = 1636.7142  # K5
$ws.Cells.Item(5, 13).Value = -1524.7142  # M5

# Row 32 (Item ID G32=4731)
$ws.Cells.Item(32, 8).Value = 50000750  # H32
$ws.Cells.Item(32, 10).Value = 50000750  # J32
$ws.Cells.Item(32, 12).Value = 150002250  # L32
$ws.Cells.Item(32, 14).Value = -150002816  # N32

# Row 36 (Item ID G36=4732)
$ws.Cells.Item(36, 8).Value = 1420.25  # H36
$ws.Cells.Item(36, 9).Value = 191  # I36
$ws.Cells.Item(36, 10).Value = 2649.5  # J36
$ws.Cells.Item(36, 11).Value = 573  # K36
$ws.Cells.Item(36, 12).Value = 7948.5  # L36
$ws.Cells.Item(36, 13).Value = -404  # M36
$ws.Cells.Item(36, 14).Value = -8286.5  # N36

# Row 38 (Item ID G38=4860)
$ws.Cells.Item(38, 8).Value = 886.375  # H38
$ws.Cells.Item(38, 9).Value = 97.666664  # I38
$ws.Cells.Item(38, 10).Value = 1359.6  # J38
$ws.Cells.Item(38, 11).Value = 292.999992  # K38
$ws.Cells.Item(38, 12).Value = 4078.8  # L38
$ws.Cells.Item(38, 13).Value = 54.00000799999998  # M38
$ws.Cells.Item(38, 14).Value = -4772.799999999999  # N38

# Row 55 (Item ID G55=4733)
$ws.Cells.Item(55, 8).Value = 1346.0625  # H55
$ws.Cells.Item(55, 9).Value = 2580  # I55
$ws.Cells.Item(55, 11).Value = 7740  # K55
$ws.Cells.Item(55, 13).Value = -7563  # M55

# Row 68 (Item ID G68=12895)
$ws.Cells.Item(68, 8).Value = 1888.6666  # H68
$ws.Cells.Item(68, 9).Value = 724  # I68
$ws.Cells.Item(68, 10).Value = 2312.182  # J68
$ws.Cells.Item(68, 11).Value = 2172  # K68
$ws.Cells.Item(68, 12).Value = 6936.545999999999  # L68
$ws.Cells.Item(68, 13).Value = -1361  # M68
$ws.Cells.Item(68, 14).Value = -8558.545999999998  # N68

# Row 71 (Item ID G71=12895)
$ws.Cells.Item(71, 8).Value = 1888.6666  # H71
$ws.Cells.Item(71, 9).Value = 724  # I71
$ws.Cells.Item(71, 10).Value = 2312.182  # J71
$ws.Cells.Item(71, 11).Value = 6516  # K71
$ws.Cells.Item(71, 12).Value = 20809.638  # L71
$ws.Cells.Item(71, 13).Value = -2460  # M71
$ws.Cells.Item(71, 14).Value = -28921.638  # N71

# Row 80 (Item ID G80=12890)
$ws.Cells.Item(80, 8).Value = 5248.75  # H80
$ws.Cells.Item(80, 10).Value = 5500  # J80
$ws.Cells.Item(80, 12).Value = 16500  # L80
$ws.Cells.Item(80, 14).Value = -18372  # N80

# Row 83 (Item ID G83=12890)
$ws.Cells.Item(83, 8).Value = 5248.75  # H83
$ws.Cells.Item(83, 10).Value = 5500  # J83
$ws.Cells.Item(83, 12).Value = 49500  # L83
$ws.Cells.Item(83, 14).Value = -58860  # N83

# Row 92 (Item ID G92=19841)
$ws.Cells.Item(92, 8).Value = 822.5  # H92
$ws.Cells.Item(92, 9).Value = 890  # I92
$ws.Cells.Item(92, 10).Value = 800  # J92
$ws.Cells.Item(92, 11).Value = 2670  # K92
$ws.Cells.Item(92, 12).Value = 2400  # L92
$ws.Cells.Item(92, 13).Value = -1422  # M92
$ws.Cells.Item(92, 14).Value = -4896  # N92

# Row 118 (Item ID G118=27872)
$ws.Cells.Item(118, 8).Value = 3710  # H118
$ws.Cells.Item(118, 9).Value = 1465.6  # I118
$ws.Cells.Item(118, 10).Value = 14932  # J118
$ws.Cells.Item(118, 11).Value = 4396.799999999999  # K118
$ws.Cells.Item(118, 12).Value = 44796  # L118
$ws.Cells.Item(118, 13).Value = -3153.799999999999  # M118
$ws.Cells.Item(118, 14).Value = -47282  # N118

# Row 128 (Item ID G128=41814)
$ws.Cells.Item(128, 8).Value = 177411.42  # H128
$ws.Cells.Item(128, 9).Value = 177411.42  # I128
$ws.Cells.Item(128, 11).Value = 532234.26  # K128
$ws.Cells.Item(128, 13).Value = -527254.26  # M128

# Row 135 (Item ID G135=43974)
$ws.Cells.Item(135, 8).Value = 645.8889  # H135
$ws.Cells.Item(135, 9).Value = 545.5714  # I135
$ws.Cells.Item(135, 11).Value = 4910.1426  # K135
$ws.Cells.Item(135, 13).Value = -2375.1426  # M135

# Row 137 (Item ID G137=44088)
$ws.Cells.Item(137, 8).Value = 12736.471  # H137
$ws.Cells.Item(137, 9).Value = 9977  # I137
$ws.Cells.Item(137, 11).Value = 29931  # K137
$ws.Cells.Item(137, 13).Value = -24831  # M137

$ws = $wb.Worksheets.Item(6)  # GSM
# Row 17 (Item ID G17=2445)
$ws.Cells.Item(17, 8).Value = 204.875  # H17
$ws.Cells.Item(17, 9).Value = 90  # I17
$ws.Cells.Item(17, 10).Value = 549.5  # J17
$ws.Cells.Item(17, 11).Value = 90  # K17
$ws.Cells.Item(17, 12).Value = 549.5  # L17
$ws.Cells.Item(17, 13).Value = 78  # M17
$ws.Cells.Item(17, 14).Value = -885.5  # N17

# Row 33 (Item ID G33=4450)
$ws.Cells.Item(33, 8).Value = 17219.666  # H33
$ws.Cells.Item(33, 9).Value = 14993.333  # I33
$ws.Cells.Item(33, 10).Value = 19446  # J33
$ws.Cells.Item(33, 11).Value = 14993.333  # K33
$ws.Cells.Item(33, 12).Value = 19446  # L33
$ws.Cells.Item(33, 13).Value = -14741.333  # M33
$ws.Cells.Item(33, 14).Value = -19950  # N33

# Row 70 (Item ID G70=14146)
$ws.Cells.Item(70, 8).Value = 6107.125  # H70
$ws.Cells.Item(70, 9).Value = 5309.8  # I70
$ws.Cells.Item(70, 11).Value = 5309.8  # K70
$ws.Cells.Item(70, 13).Value = -5039.8  # M70

# Row 73 (Item ID G73=14146)
$ws.Cells.Item(73, 8).Value = 6107.125  # H73
$ws.Cells.Item(73, 9).Value = 5309.8  # I73
$ws.Cells.Item(73, 11).Value = 5309.8  # K73
$ws.Cells.Item(73, 13).Value = -4373.8  # M73

# Row 80 (Item ID G80=12521)
$ws.Cells.Item(80, 8).Value = 50717280  # H80
$ws.Cells.Item(80, 9).Value = 64547736  # I80
$ws.Cells.Item(80, 11).Value = 64547736  # K80
$ws.Cells.Item(80, 13).Value = -64546738  # M80

# Row 83 (Item ID G83=12521)
$ws.Cells.Item(83, 8).Value = 50717280  # H83
$ws.Cells.Item(83, 9).Value = 64547736  # I83
$ws.Cells.Item(83, 11).Value = 322738680  # K83
$ws.Cells.Item(83, 13).Value = -322733688  # M83

# Row 102 (Item ID G102=36169)
$ws.Cells.Item(102, 8).Value = 9388  # H102
$ws.Cells.Item(102, 9).Value = 6448.6924  # I102
$ws.Cells.Item(102, 11).Value = 6448.6924  # K102
$ws.Cells.Item(102, 13).Value = -4826.6924  # M102

# Row 113 (Item ID G113=27710)
$ws.Cells.Item(113, 8).Value = 2913.1667  # H113
$ws.Cells.Item(113, 9).Value = 2331.2222  # I113
$ws.Cells.Item(113, 10).Value = 4659  # J113
$ws.Cells.Item(113, 11).Value = 2331.2222  # K113
$ws.Cells.Item(113, 12).Value = 4659  # L113
$ws.Cells.Item(113, 13).Value = -161.2222000000002  # M113
$ws.Cells.Item(113, 14).Value = -8999  # N113

# Row 122 (Item ID G122=36182)
$ws.Cells.Item(122, 8).Value = 3573.3667  # H122
$ws.Cells.Item(122, 9).Value = 3086.6428  # I122
$ws.Cells.Item(122, 11).Value = 9259.928400000001  # K122
$ws.Cells.Item(122, 13).Value = -6809.928400000001  # M122

# Row 126 (Item ID G126=36184)
$ws.Cells.Item(126, 8).Value = 6176.5713  # H126
$ws.Cells.Item(126, 9).Value = 5809.1113  # I126
$ws.Cells.Item(126, 10).Value = 6838  # J126
$ws.Cells.Item(126, 11).Value = 17427.3339  # K126
$ws.Cells.Item(126, 12).Value = 20514  # L126
$ws.Cells.Item(126, 13).Value = -14957.3339  # M126
$ws.Cells.Item(126, 14).Value = -25454  # N126

# Row 132 (Item ID G132=44008)
$ws.Cells.Item(132, 8).Value = 7153  # H132
$ws.Cells.Item(132, 10).Value = 7183.6  # J132
$ws.Cells.Item(132, 12).Value = 21550.8  # L132
$ws.Cells.Item(132, 14).Value = -26610.8  # N132

$ws = $wb.Worksheets.Item(7)  # LTW
# Row 46 (Item ID G46=5282)
$ws.Cells.Item(46, 8).Value = 3150.1365  # H46
$ws.Cells.Item(46, 9).Value = 2646.3635  # I46
$ws.Cells.Item(46, 10).Value = 3653.9092  # J46
$ws.Cells.Item(46, 11).Value = 2646.3635  # K46
$ws.Cells.Item(46, 12).Value = 3653.9092  # L46
$ws.Cells.Item(46, 13).Value = -2458.3635  # M46
$ws.Cells.Item(46, 14).Value = -4029.9092  # N46

# Row 55 (Item ID G55=5284)
$ws.Cells.Item(55, 8).Value = 1303.5  # H55
$ws.Cells.Item(55, 9).Value = 1190.3334  # I55
$ws.Cells.Item(55, 10).Value = 1575.1  # J55
$ws.Cells.Item(55, 11).Value = 1190.3334  # K55
$ws.Cells.Item(55, 12).Value = 1575.1  # L55
$ws.Cells.Item(55, 13).Value = -1017.3334  # M55
$ws.Cells.Item(55, 14).Value = -1921.1  # N55

# Row 82 (Item ID G82=12565)
$ws.Cells.Item(82, 8).Value = 423.24  # H82
$ws.Cells.Item(82, 9).Value = 409.62888  # I82
$ws.Cells.Item(82, 10).Value = 863.3333  # J82
$ws.Cells.Item(82, 11).Value = 409.62888  # K82
$ws.Cells.Item(82, 12).Value = 863.3333  # L82
$ws.Cells.Item(82, 13).Value = -48.62887999999998  # M82
$ws.Cells.Item(82, 14).Value = -1585.3333  # N82

# Row 85 (Item ID G85=12565)
$ws.Cells.Item(85, 8).Value = 423.24  # H85
$ws.Cells.Item(85, 9).Value = 409.62888  # I85
$ws.Cells.Item(85, 10).Value = 863.3333  # J85
$ws.Cells.Item(85, 11).Value = 409.62888  # K85
$ws.Cells.Item(85, 12).Value = 863.3333  # L85
$ws.Cells.Item(85, 13).Value = 838.37112  # M85
$ws.Cells.Item(85, 14).Value = -3359.3333  # N85

# Row 93 (Item ID G93=19993)
$ws.Cells.Item(93, 8).Value = 11111845  # H93
$ws.Cells.Item(93, 9).Value = 20000674  # I93
$ws.Cells.Item(93, 10).Value = 808.125  # J93
$ws.Cells.Item(93, 11).Value = 20000674  # K93
$ws.Cells.Item(93, 12).Value = 808.125  # L93
$ws.Cells.Item(93, 13).Value = -19999426  # M93
$ws.Cells.Item(93, 14).Value = -3304.125  # N93

# Row 100 (Item ID G100=19995)
$ws.Cells.Item(100, 8).Value = 83336360  # H100
$ws.Cells.Item(100, 9).Value = 250003410  # I100
$ws.Cells.Item(100, 10).Value = 2837.5  # J100
$ws.Cells.Item(100, 11).Value = 250003410  # K100
$ws.Cells.Item(100, 12).Value = 2837.5  # L100
$ws.Cells.Item(100, 13).Value = -250002869  # M100
$ws.Cells.Item(100, 14).Value = -3919.5  # N100

# Row 122 (Item ID G122=36247)
$ws.Cells.Item(122, 8).Value = 0  # H122
$ws.Cells.Item(122, 9).Value = 0  # I122
$ws.Cells.Item(122, 11).Value = 0  # K122
$ws.Cells.Item(122, 13).ClearContents()  # M122

# Row 132 (Item ID G132=44058)
$ws.Cells.Item(132, 8).Value = 47169.523  # H132
$ws.Cells.Item(132, 9).Value = 51185.668  # I132
$ws.Cells.Item(132, 10).Value = 5000  # J132
$ws.Cells.Item(132, 11).Value = 153557.004  # K132
$ws.Cells.Item(132, 12).Value = 15000  # L132
$ws.Cells.Item(132, 13).Value = -151027.004  # M132
$ws.Cells.Item(132, 14).Value = -20060  # N132

# Row 136 (Item ID G136=44060)
$ws.Cells.Item(136, 8).Value = 5822604.5  # H136
$ws.Cells.Item(136, 9).Value = 7841641  # I136
$ws.Cells.Item(136, 10).Value = 17873.875  # J136
$ws.Cells.Item(136, 11).Value = 23524923  # K136
$ws.Cells.Item(136, 12).Value = 53621.625  # L136
$ws.Cells.Item(136, 13).Value = -23522373  # M136
$ws.Cells.Item(136, 14).Value = -58721.625  # N136

$ws = $wb.Worksheets.Item(8)  # WVR
# Row 99 (Item ID G99=19640)
$ws.Cells.Item(99, 8).Value = 44526  # H99
$ws.Cells.Item(99, 9).Value = 44526  # I99
$ws.Cells.Item(99, 11).Value = 44526  # K99
$ws.Cells.Item(99, 13).Value = -41531  # M99

# Row 100 (Item ID G100=19981)
$ws.Cells.Item(100, 8).Value = 590.73334  # H100
$ws.Cells.Item(100, 9).Value = 459.36365  # I100
$ws.Cells.Item(100, 10).Value = 952  # J100
$ws.Cells.Item(100, 11).Value = 918.7273  # K100
$ws.Cells.Item(100, 12).Value = 1904  # L100
$ws.Cells.Item(100, 13).Value = -377.7273  # M100
$ws.Cells.Item(100, 14).Value = -2986  # N100

# Row 113 (Item ID G113=27752)
$ws.Cells.Item(113, 8).Value = 951.7778  # H113
$ws.Cells.Item(113, 9).Value = 492.2857  # I113
$ws.Cells.Item(113, 10).Value = 2560  # J113
$ws.Cells.Item(113, 11).Value = 1476.8571  # K113
$ws.Cells.Item(113, 12).Value = 7680  # L113
$ws.Cells.Item(113, 13).Value = 693.1428999999998  # M113
$ws.Cells.Item(113, 14).Value = -12020  # N113

# Row 122 (Item ID G122=36208)
$ws.Cells.Item(122, 8).Value = 13231.611  # H122
$ws.Cells.Item(122, 9).Value = 9082.615  # I122
$ws.Cells.Item(122, 10).Value = 24019  # J122
$ws.Cells.Item(122, 11).Value = 27247.845  # K122
$ws.Cells.Item(122, 12).Value = 72057  # L122
$ws.Cells.Item(122, 13).Value = -24797.845  # M122
$ws.Cells.Item(122, 14).Value = -76957  # N122

# Row 132 (Item ID G132=44029)
$ws.Cells.Item(132, 8).Value = 5802.3125  # H132
$ws.Cells.Item(132, 9).Value = 5889.1333  # I132
$ws.Cells.Item(132, 10).Value = 4500  # J132
$ws.Cells.Item(132, 11).Value = 17667.3999  # K132
$ws.Cells.Item(132, 12).Value = 13500  # L132
$ws.Cells.Item(132, 13).Value = -15137.3999  # M132
$ws.Cells.Item(132, 14).Value = -18560  # N132

# Row 136 (Item ID G136=44031)
$ws.Cells.Item(136, 8).Value = 961.75  # H136
$ws.Cells.Item(136, 9).Value = 961.75  # I136
$ws.Cells.Item(136, 11).Value = 2885.25  # K136
$ws.Cells.Item(136, 13).Value = -335.25  # M136
